{"js": "// Remove the inline picture (\"photo\") from the document.\n//\n// The source document has exactly one inline picture, sitting as the\n// first run of the first (centered) paragraph that also holds the\n// \"Offre d'emploi\" heading text. Word keeps a hidden \"_GoBack\" bookmark\n// that tracks the location of the last edit; after deleting the\n// picture, Word re-stamps that bookmark at the edit point (collapsing\n// it where the picture used to be) and removes it from wherever it\n// previously sat. We reproduce both effects explicitly so the saved\n// OOXML matches exactly what Word itself produces for this edit.\n\nconst body = context.document.body;\n\n// Locate the paragraph that contains the picture and remember a\n// collapsed range at its very start - this is where the picture (and\n// therefore the relocated \"_GoBack\" bookmark) lives.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nconst pictures = body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nif (pictures.items.length > 0) {\n  const picture = pictures.items[0];\n\n  // `InlinePicture.paragraph` isn't reliable here, so find the owning\n  // paragraph by checking each paragraph's own inlinePictures instead.\n  for (const p of paragraphs.items) {\n    p.inlinePictures.load(\"items\");\n  }\n  await context.sync();\n\n  let picParagraph = null;\n  for (const p of paragraphs.items) {\n    if (p.inlinePictures.items.length > 0) {\n      picParagraph = p;\n      break;\n    }\n  }\n\n  const picStart = (picParagraph || paragraphs.items[0]).getRange(\"Start\");\n\n  // Drop the old \"_GoBack\" bookmark (Word only ever keeps a single\n  // instance of it, wherever the last edit happened).\n  context.document.deleteBookmark(\"_GoBack\");\n\n  // Delete the picture itself - this is the actual user-visible edit.\n  picture.delete();\n\n  // Re-create \"_GoBack\" collapsed at the spot the picture used to\n  // occupy, matching Word's own last-edit bookmark bookkeeping.\n  picStart.insertBookmark(\"_GoBack\");\n\n  await context.sync();\n}\n", "ps1": "# Remove the inline picture (\"photo\") from the document.\n#\n# The source document has exactly one inline picture, sitting as the\n# first run of the first (centered) paragraph that also holds the\n# \"Offre d'emploi\" heading text. Word keeps a hidden \"_GoBack\" bookmark\n# that tracks the location of the last edit; after deleting the\n# picture, Word re-stamps that bookmark at the edit point (collapsing\n# it where the picture used to be) and removes it from wherever it\n# previously sat. We reproduce both effects explicitly so the saved\n# OOXML matches exactly what Word itself produces for this edit.\n\n$d = $word.ActiveDocument\n\nif ($d.InlineShapes.Count -gt 0) {\n    $shp = $d.InlineShapes(1)\n\n    # Remember a collapsed range at the very start of the picture -\n    # this is where it (and therefore the relocated \"_GoBack\" bookmark)\n    # lives. Duplicate so later edits don't shift this particular range\n    # object out from under us.\n    $picRange = $shp.Range.Duplicate\n    $picRange.Collapse(1)  # wdCollapseStart\n\n    # Drop the old \"_GoBack\" bookmark (Word only ever keeps a single\n    # instance of it, wherever the last edit happened).\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks.Item(\"_GoBack\").Delete()\n    }\n\n    # Delete the picture itself - this is the actual user-visible edit.\n    $shp.Delete()\n\n    # Re-create \"_GoBack\" collapsed at the spot the picture used to\n    # occupy, matching Word's own last-edit bookmark bookkeeping.\n    $d.Bookmarks.Add(\"_GoBack\", $picRange)\n}\n"}
